$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.064.87'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '2.360.58'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('E5').Value = '  +1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '239.80'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.37'
$ws.Range('E7').Value = '  +1.06%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.598'
$ws.Range('E9').Value = '  +8.81%  '
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.29'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '32.43'
$ws.Range('E12').Value = '  +8.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.30'
$ws.Range('E13').Value = '  +7.66%  '
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').Value = '2.710.53'
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '16.60'
$ws.Range('E16').Value = '  -2.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.904'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = '2.357.96'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = '43.912.86'
$ws.Range('E19').Value = '  -1.42%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.92'
$ws.Range('E20').Value = '  +6.60%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000102'
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '77.18'
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '258.86'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.97'
$ws.Range('E24').Value = '  +22.97%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('B26').Value = 'WEMIXToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.66'
$ws.Range('E26').Value = '  -2.88%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.49'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.85'
$ws.Range('E28').Value = '  +3.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.28'
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.79'
$ws.Range('E30').Value = '  +1.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '176.09'
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('E33').Value = '  +2.24%  '
$ws.Range('E34').Value = '  +2.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.25'
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.54'
$ws.Range('E36').Value = '  +6.03%  '
$ws.Range('E37').Value = '  -4.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.34'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.36'
$ws.Range('E39').Value = '  -3.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.112'
$ws.Range('E41').Value = '  +14.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.204'
$ws.Range('E42').Value = '  +10.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.06'
$ws.Range('E43').Value = '  -6.15%  '
$ws.Range('E44').Value = '  +1.06%  '
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.70'
$ws.Range('E46').Value = '  +5.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '58.23'
$ws.Range('E47').Value = '  +10.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.50'
$ws.Range('E48').Value = '  +5.72%  '
$ws.Range('E49').Value = '  -1.11%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '100.87'
$ws.Range('E50').Value = '  +2.07%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.18'
$ws.Range('E51').Value = '  +0.02%  '
